$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(3, 6).Value = 6494
$ws.Cells.Item(4, 6).Value = 1039
$ws.Cells.Item(5, 6).Value = 647
$ws.Cells.Item(6, 6).Value = 1440
$ws.Cells.Item(7, 6).Value = 3190
$ws.Cells.Item(9, 6).Value = 564
$ws.Cells.Item(10, 6).Value = 2110
$ws.Cells.Item(11, 6).Value = 451
$ws.Cells.Item(12, 6).Value = 381
$ws.Cells.Item(13, 6).Value = 222
$ws.Cells.Item(14, 6).Value = 110
$ws.Cells.Item(15, 6).Value = 246
$ws.Cells.Item(16, 6).Value = 1039
$ws.Cells.Item(17, 6).Value = 407
$ws.Cells.Item(18, 6).Value = 66
$ws.Cells.Item(19, 6).Value = 164
$ws.Cells.Item(20, 6).Value = 4105
$ws.Cells.Item(21, 6).Value = 1232
$ws.Cells.Item(22, 6).Value = 3210
$ws.Cells.Item(23, 6).Value = 313
$ws.Cells.Item(24, 6).Value = 107
$ws.Cells.Item(25, 6).Value = 2999
$ws.Cells.Item(26, 6).Value = 4667
$ws.Cells.Item(27, 6).Value = 119
$ws.Cells.Item(29, 6).Value = 511
$ws.Cells.Item(30, 6).Value = 3046
$ws.Cells.Item(31, 6).Value = 299
$ws.Cells.Item(32, 6).Value = 44
$ws.Cells.Item(33, 6).Value = 112
$ws.Cells.Item(34, 6).Value = 70
$ws.Cells.Item(36, 6).Value = 1103
$ws.Cells.Item(37, 6).Value = 1352
$ws.Cells.Item(39, 6).Value = 1234
$ws.Cells.Item(40, 6).Value = 798
$ws.Cells.Item(41, 6).Value = 7
$ws.Cells.Item(42, 6).Value = 733
$ws.Cells.Item(43, 6).Value = 479
$ws.Cells.Item(44, 6).Value = 42
$ws.Cells.Item(45, 6).Value = 209
$ws.Cells.Item(46, 6).Value = 43
$ws.Cells.Item(47, 6).Value = 88
$ws.Cells.Item(48, 6).Value = 350
$ws.Cells.Item(49, 6).Value = 3673

$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(7, 6).Value = 966
$ws.Cells.Item(20, 6).Value = 46

$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2, 6).Value = 1788

$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(2, 6).Value = 6494
$ws.Cells.Item(4, 6).Value = 1789
$ws.Cells.Item(5, 6).Value = 647
$ws.Cells.Item(6, 6).Value = 1440
$ws.Cells.Item(7, 6).Value = 3190
$ws.Cells.Item(9, 6).Value = 2110
$ws.Cells.Item(10, 6).Value = 451
$ws.Cells.Item(11, 6).Value = 381
$ws.Cells.Item(13, 6).Value = 222
$ws.Cells.Item(14, 6).Value = 966
$ws.Cells.Item(16, 6).Value = 110
$ws.Cells.Item(17, 6).Value = 246
$ws.Cells.Item(18, 6).Value = 1039
$ws.Cells.Item(20, 6).Value = 407
$ws.Cells.Item(21, 6).Value = 164
$ws.Cells.Item(22, 6).Value = 4105
$ws.Cells.Item(24, 6).Value = 1232
$ws.Cells.Item(26, 6).Value = 3210
$ws.Cells.Item(27, 6).Value = 2999
$ws.Cells.Item(28, 6).Value = 4667
$ws.Cells.Item(30, 6).Value = 3046
$ws.Cells.Item(31, 6).Value = 299
$ws.Cells.Item(32, 6).Value = 44
$ws.Cells.Item(33, 6).Value = 112
$ws.Cells.Item(35, 6).Value = 1103
$ws.Cells.Item(36, 6).Value = 1352
$ws.Cells.Item(38, 6).Value = 1234
$ws.Cells.Item(39, 6).Value = 798
$ws.Cells.Item(41, 6).Value = 479
$ws.Cells.Item(42, 6).Value = 46
$ws.Cells.Item(43, 6).Value = 42
$ws.Cells.Item(45, 6).Value = 209
$ws.Cells.Item(46, 6).Value = 43
$ws.Cells.Item(47, 6).Value = 88
$ws.Cells.Item(48, 6).Value = 350
$ws.Cells.Item(49, 6).Value = 3673
